$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046958558114491
$ws.Range("D2").Value = 1.054052416041449
$ws.Range("E2").Value = 1.054185947263162
$ws.Range("F2").Value = 1.06420249081073
$ws.Range("I2").Value = 1.044111746295256
$ws.Range("J2").Value = 1.052010179334492
$ws.Range("K2").Value = 1.056796689637625
$ws.Range("L2").Value = 1.056929852692584
$ws.Range("M2").Value = 1.06691906415494
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047830477798015
$ws.Range("D3").Value = 1.054745458157885
$ws.Range("E3").Value = 1.054953632148761
$ws.Range("F3").Value = 1.065032014408452
$ws.Range("I3").Value = 1.0443166027211
$ws.Range("J3").Value = 1.052530759808397
$ws.Range("K3").Value = 1.057303205803162
$ws.Range("L3").Value = 1.057510847156796
$ws.Range("M3").Value = 1.067563711040477
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048395382249962
$ws.Range("D4").Value = 1.055194523577733
$ws.Range("E4").Value = 1.055451374365455
$ws.Range("F4").Value = 1.065569848291315
$ws.Range("I4").Value = 1.044448320777421
$ws.Range("J4").Value = 1.052867653377263
$ws.Range("K4").Value = 1.057630897147674
$ws.Range("L4").Value = 1.057887123717556
$ws.Range("M4").Value = 1.067981264146262
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048633037433901
$ws.Range("D5").Value = 1.055383457526973
$ws.Range("E5").Value = 1.055660862469504
$ws.Range("F5").Value = 1.065796209485627
$ws.Range("I5").Value = 1.044503493774262
$ws.Range("J5").Value = 1.053009292477182
$ws.Range("K5").Value = 1.057768643242555
$ws.Range("L5").Value = 1.058045388973859
$ws.Range("M5").Value = 1.068156903263886
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048672950669266
$ws.Range("D6").Value = 1.055415188927776
$ws.Range("E6").Value = 1.055696050304754
$ws.Range("F6").Value = 1.065834231472323
$ws.Range("I6").Value = 1.044512745739651
$ws.Range("J6").Value = 1.053033074797942
$ws.Range("K6").Value = 1.057791770478991
$ws.Range("L6").Value = 1.058071966973496
$ws.Range("M6").Value = 1.068186399650323
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048398557146728
$ws.Range("D7").Value = 1.055197047547639
$ws.Range("E7").Value = 1.055454172625719
$ws.Range("F7").Value = 1.065572871937847
$ws.Range("I7").Value = 1.044449058793465
$ws.Range("J7").Value = 1.052869545931937
$ws.Range("K7").Value = 1.057632737779722
$ws.Range("L7").Value = 1.057889238158657
$ws.Range("M7").Value = 1.067983610654623
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047253078624127
$ws.Range("D8").Value = 1.054286503252453
$ws.Range("E8").Value = 1.054445182064364
$ws.Range("F8").Value = 1.06448260811233
$ws.Range("I8").Value = 1.044181151514525
$ws.Range("J8").Value = 1.052186102219871
$ws.Range("K8").Value = 1.05696788044788
$ws.Range("L8").Value = 1.057126132408524
$ws.Range("M8").Value = 1.067136836953073
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045240139592754
$ws.Range("D9").Value = 1.052686841418247
$ws.Range("E9").Value = 1.05267494202276
$ws.Range("F9").Value = 1.062569754781738
$ws.Range("I9").Value = 1.043702680197919
$ws.Range("J9").Value = 1.0509821752574
$ws.Range("K9").Value = 1.055795928480066
$ws.Range("L9").Value = 1.055784066568081
$ws.Range("M9").Value = 1.065648026917974
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043902001160398
$ws.Range("D10").Value = 1.051623756307175
$ws.Range("E10").Value = 1.051500088601749
$ws.Range("F10").Value = 1.061300235135973
$ws.Range("I10").Value = 1.043379449530866
$ws.Range("J10").Value = 1.050179895633573
$ws.Range("K10").Value = 1.055014446313408
$ws.Range("L10").Value = 1.054891206176154
$ws.Range("M10").Value = 1.064657809341672
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04332349885199
$ws.Range("D11").Value = 1.051164247211499
$ws.Range("E11").Value = 1.050992645845067
$ws.Range("F11").Value = 1.060751899096817
$ws.Range("I11").Value = 1.043238488100151
$ws.Range("J11").Value = 1.049832596262965
$ws.Range("K11").Value = 1.054676029952137
$ws.Range("L11").Value = 1.054505044747059
$ws.Range("M11").Value = 1.064229605046871
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043108757238664
$ws.Range("D12").Value = 1.050993689153792
$ws.Range("E12").Value = 1.050804352669001
$ws.Range("F12").Value = 1.060548431123664
$ws.Range("I12").Value = 1.043185979158959
$ws.Range("J12").Value = 1.049703609095772
$ws.Range("K12").Value = 1.054550323988236
$ws.Range("L12").Value = 1.054361676722299
$ws.Range("M12").Value = 1.064070637623058
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043154813686889
$ws.Range("D13").Value = 1.051030268792048
$ws.Range("E13").Value = 1.050844733400948
$ws.Range("F13").Value = 1.060592066230378
$ws.Range("I13").Value = 1.043197249262874
$ws.Range("J13").Value = 1.049731276568377
$ws.Range("K13").Value = 1.054577288460783
$ws.Range("L13").Value = 1.054392426480678
$ws.Range("M13").Value = 1.064104732725039
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043305745379946
$ws.Range("D14").Value = 1.051150146287525
$ws.Range("E14").Value = 1.050977077511004
$ws.Range("F14").Value = 1.060735076107066
$ws.Range("I14").Value = 1.043234150742145
$ws.Range("J14").Value = 1.049821933824146
$ws.Range("K14").Value = 1.054665639112478
$ws.Range("L14").Value = 1.054493192474676
$ws.Range("M14").Value = 1.064216462971689
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04339875796408
$ws.Range("D15").Value = 1.051224023265178
$ws.Range("E15").Value = 1.05105864480267
$ws.Range("F15").Value = 1.060823216892016
$ws.Range("I15").Value = 1.043256867160706
$ws.Range("J15").Value = 1.049877792815791
$ws.Range("K15").Value = 1.054720074500278
$ws.Range("L15").Value = 1.054555286985651
$ws.Range("M15").Value = 1.064285315203361
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043940413984006
$ws.Range("D16").Value = 1.051654269713123
$ws.Range("E16").Value = 1.051533792961221
$ws.Range("F16").Value = 1.061336655516735
$ws.Range("I16").Value = 1.043388783659832
$ws.Range("J16").Value = 1.050202946820233
$ws.Range("K16").Value = 1.05503690538899
$ws.Range("L16").Value = 1.05491684410781
$ws.Range("M16").Value = 1.064686239922704
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044280428066073
$ws.Range("D17").Value = 1.051924371098065
$ws.Range("E17").Value = 1.051832183933515
$ws.Range("F17").Value = 1.061659091340694
$ws.Range("I17").Value = 1.043471264055973
$ws.Range("J17").Value = 1.050406933160311
$ws.Range("K17").Value = 1.055235638117655
$ws.Range("L17").Value = 1.055143761655923
$ws.Range("M17").Value = 1.064937882201134
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044478841230515
$ws.Range("D18").Value = 1.052081995061418
$ws.Range("E18").Value = 1.052006353380045
$ws.Range("F18").Value = 1.061847295174108
$ws.Range("I18").Value = 1.043519276901784
$ws.Range("J18").Value = 1.050525923826597
$ws.Range("K18").Value = 1.05535155262544
$ws.Range("L18").Value = 1.055276162433554
$ws.Range("M18").Value = 1.065084715430633
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044546509988699
$ws.Range("D19").Value = 1.052135754001368
$ws.Range("E19").Value = 1.052065761473198
$ws.Range("F19").Value = 1.061911490227476
$ws.Range("I19").Value = 1.043535631639219
$ws.Range("J19").Value = 1.05056649803606
$ws.Range("K19").Value = 1.055391075938611
$ws.Range("L19").Value = 1.05532131502983
$ws.Range("M19").Value = 1.065134790970212
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044243938586328
$ws.Range("D20").Value = 1.051895383653755
$ws.Range("E20").Value = 1.051800156665931
$ws.Range("F20").Value = 1.061624483311849
$ws.Range("I20").Value = 1.043462424674001
$ws.Range("J20").Value = 1.050385046428027
$ws.Range("K20").Value = 1.055214316269955
$ws.Range("L20").Value = 1.055119411029089
$ws.Range("M20").Value = 1.064910877712589
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043261295883267
$ws.Range("D21").Value = 1.051114841913149
$ws.Range("E21").Value = 1.05093810010824
$ws.Range("F21").Value = 1.060692957492874
$ws.Range("I21").Value = 1.043223288299918
$ws.Range("J21").Value = 1.049795237089495
$ws.Range("K21").Value = 1.054639622121151
$ws.Range("L21").Value = 1.054463517470015
$ws.Range("M21").Value = 1.064183558792073
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04264427964534
$ws.Range("D22").Value = 1.050624803130376
$ws.Range("E22").Value = 1.050397212305951
$ws.Range("F22").Value = 1.060108477218626
$ws.Range("I22").Value = 1.043072068645616
$ws.Range("J22").Value = 1.04942448961842
$ws.Range("K22").Value = 1.054278272074467
$ws.Range("L22").Value = 1.054051534125433
$ws.Range("M22").Value = 1.06372676703718
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04297129428722
$ws.Range("D23").Value = 1.050884513220619
$ws.Range("E23").Value = 1.050683840245416
$ws.Range("F23").Value = 1.060418206206289
$ws.Range("I23").Value = 1.043152314810954
$ws.Range("J23").Value = 1.049621020987052
$ws.Range("K23").Value = 1.054469831730156
$ws.Range("L23").Value = 1.05426989560576
$ws.Range("M23").Value = 1.063968872817207
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044260426337666
$ws.Range("D24").Value = 1.051908481590094
$ws.Range("E24").Value = 1.05181462802962
$ws.Range("F24").Value = 1.061640120781333
$ws.Range("I24").Value = 1.04346641910871
$ws.Range("J24").Value = 1.05039493607295
$ws.Range("K24").Value = 1.055223950704366
$ws.Range("L24").Value = 1.055130413894628
$ws.Range("M24").Value = 1.064923079710233
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045759866150321
$ws.Range("D25").Value = 1.053099808588945
$ws.Range("E25").Value = 1.053131664175613
$ws.Range("F25").Value = 1.063063274151133
$ws.Range("I25").Value = 1.043827128696101
$ws.Range("J25").Value = 1.051293365266783
$ws.Range("K25").Value = 1.056098943913724
$ws.Range("L25").Value = 1.056130702946362
$ws.Range("M25").Value = 1.066032517586609
